$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "42.478.49"
Set-TextValue $ws.Range("E2") "  -1.36%  "
Set-TextValue $ws.Range("D3") "2.183.81"
Set-TextValue $ws.Range("E3") "  -2.38%  "
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "251.65"
Set-TextValue $ws.Range("E5") "  +2.32%  "
Set-TextValue $ws.Range("E6") "  -0.91%  "
Set-TextValue $ws.Range("D7") "75.07"
Set-TextValue $ws.Range("E7") "  -0.82%  "
Set-TextValue $ws.Range("E8") "  +0.04%  "
Set-TextValue $ws.Range("D9") "0.582"
Set-TextValue $ws.Range("E9") "  -5.62%  "
Set-TextValue $ws.Range("D10") "40.25"
Set-TextValue $ws.Range("E10") "  -2.15%  "
Set-TextValue $ws.Range("D11") "0.0908"
Set-TextValue $ws.Range("E11") "  -2.24%  "
Set-TextValue $ws.Range("D12") "0.101"
Set-TextValue $ws.Range("E12") "  -0.57%  "
Set-TextValue $ws.Range("E13") "  -2.95%  "
Set-TextValue $ws.Range("D14") "2.512.83"
Set-TextValue $ws.Range("E14") "  -2.18%  "
Set-TextValue $ws.Range("D15") "14.14"
Set-TextValue $ws.Range("E15") "  -3.88%  "
Set-TextValue $ws.Range("D16") "2.184.01"
Set-TextValue $ws.Range("E16") "  -2.07%  "
Set-TextValue $ws.Range("D17") "0.766"
Set-TextValue $ws.Range("E17") "  -5.74%  "
Set-TextValue $ws.Range("D18") "42.388.20"
Set-TextValue $ws.Range("E18") "  -1.29%  "
Set-TextValue $ws.Range("E19") "  -3.37%  "
Set-TextValue $ws.Range("D20") "70.68"
Set-TextValue $ws.Range("E20") "  -0.63%  "
Set-TextValue $ws.Range("E21") "  -2.42%  "
Set-TextValue $ws.Range("D22") "226.74"
Set-TextValue $ws.Range("E22") "  -1.62%  "
Set-TextValue $ws.Range("D23") "9.43"
Set-TextValue $ws.Range("E23") "  -9.47%  "
Set-TextValue $ws.Range("E24") "  -4.24%  "
Set-TextValue $ws.Range("E25") "  +0.04%  "
Set-TextValue $ws.Range("D26") "10.43"
Set-TextValue $ws.Range("E26") "  -5.11%  "
Set-TextValue $ws.Range("D27") "3.39"
Set-TextValue $ws.Range("E27") "  +0.85%  "
Set-TextValue $ws.Range("D28") "2.20"
Set-TextValue $ws.Range("E28") "  +4.36%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D29") "37.89"
Set-TextValue $ws.Range("E29") "  +1.19%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D30") "2.15"
Set-TextValue $ws.Range("E30") "  -4.33%  "
Set-TextValue $ws.Range("D31") "172.43"
Set-TextValue $ws.Range("E31") "  -0.94%  "
Set-TextValue $ws.Range("E32") "  -1.70%  "
Set-TextValue $ws.Range("D33") "0.0820"
Set-TextValue $ws.Range("E33") "  +3.27%  "
Set-TextValue $ws.Range("D34") "5.12"
Set-TextValue $ws.Range("E34") "  -4.62%  "
Set-TextValue $ws.Range("E35") "  -1.76%  "
Set-TextValue $ws.Range("E36") "  -4.15%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D37") "0.0337"
Set-TextValue $ws.Range("E37") "  +1.22%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D38") "4.21"
Set-TextValue $ws.Range("E38") "  -3.61%  "
Set-TextValue $ws.Range("D39") "11.91"
Set-TextValue $ws.Range("E39") "  -9.43%  "
Set-TextValue $ws.Range("E40") "  -3.77%  "
Set-TextValue $ws.Range("D41") "2.60"
Set-TextValue $ws.Range("E41") "  +11.63%  "
Set-TextValue $ws.Range("D42") "5.15"
Set-TextValue $ws.Range("E42") "  -7.81%  "
Set-TextValue $ws.Range("D43") "58.75"
Set-TextValue $ws.Range("E43") "  -2.25%  "
Set-TextValue $ws.Range("E44") "  -3.40%  "
Set-TextValue $ws.Range("D45") "101.24"
Set-TextValue $ws.Range("E45") "  -3.61%  "
Set-TextValue $ws.Range("E46") "  -2.43%  "
Set-TextValue $ws.Range("D47") "0.458"
Set-TextValue $ws.Range("E47") "  +3.31%  "
Set-TextValue $ws.Range("E48") "  -4.92%  "
Set-TextValue $ws.Range("E49") "  -1.70%  "
Set-TextValue $ws.Range("E50") "  -2.40%  "
Set-TextValue $ws.Range("E51") "  -0.91%  "
